$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADI-burndown")
$ws.Range("M5:M19").Value = 0
$ws.Range("M20").Value = 3
$ws.Range("M21").Value = 3
$ws.Range("M22").Value = 0
$ws.Range("M23").Value = 3
$ws.Range("M24").Value = 35
$ws.Range("M25").Value = 15
$ws.Range("M26").Value = 10
$ws.Range("M29").Formula = "=SUM(M5:M27)"
